$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (American Samoa): fill in the "Land oversat" (translated country) column,
# matching the pattern used by other rows with no real translation (same text as column A).
$ws.Range("D3").Value = "American Samoa"

# Row 29 (Trinidad and Tobago): the "Land oversat" column previously held a
# differently-capitalised duplicate ("Trinidad And Tobago"). Point it at the
# same text used in column A ("Trinidad and Tobago") instead, which makes the
# old duplicate string unused and lets it drop out of the workbook.
$ws.Range("D29").Value = "Trinidad and Tobago"

# Row 20 (Qatar): update the matched security name list - drop the QPETRO note
# that is no longer matched.
$ws.Range("F20").Value = "['3,25 QATAR 02-06-2026 (REGS)', '4,5 QATAR 23-04-2028 (REGS)', '5.10% Qatar Government International Bond 2048', 'STATE OF QATAR 5.103% 23.04.2048', 'STATE OF QATAR 4.817% 14.03.2049', 'STATE OF QATAR 3.4% 16.04.2025', '3,75 QATAR 16-04-2030 (REGS)', '4,4 QATAR 16-04-2050 (REGS)', 'STATE OF QATAR 4.4% 16.04.2050', 'QATAR ENERGY 2.25% 12.07.2031']"

# Row 20 (Qatar): reorder the matched issuer list entries.
$ws.Range("G20").Value = "['State of Qatar', '5.10% Qatar Government International Bond 2048', 'Qatar Government International Bond', 'Qatar', 'STATE OF QATAR  ', 'QATAR ENERGY  ', 'QatarEnergy']"
